# "add real time notification" -- adds a new "Master Content Code" header
# column (F) to Sheet1, matching the look of the existing header cells,
# and moves the active selection onto the newly added header cell, as
# happens when a user types a new column header next to the last one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text in the first empty column (F) of the header row.
$ws.Range("F1").Value = "Master Content Code"

# Match the header formatting (fill/border/font/alignment) used by the
# other header cells (A1:E1) by copying the format from E1 onto F1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Give the new column a comfortable width, similar to its neighbours.
$ws.Columns.Item(6).ColumnWidth = 21.83

# Move the selection onto the newly added header cell.
$ws.Range("F1").Select() | Out-Null
